# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The data table (rows 16-27) lists, for each (worker, periodo-mora) pair,
# a row with Tipo Doc / N Doc / Nombre / Periodo Mora / Valor Mora / Salario
# Basico. The edit re-sorts those 12 rows so they are grouped by period
# (ascending 1704..1709) instead of by worker, and refreshes "Salario
# Basico" (column G) from 738000 to 781242 for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ccDoc      = "CC"
$doc1       = "1050953042"
$name1      = "JHONNY JAVIER CARRILLO GARCIA"
$doc2       = "1047378133"
$name2      = "MAURICIO JAVIER TORRES ANGULO"

$valorMora    = 29520
$salarioBasico = 781242

$periodos = @("1704", "1705", "1706", "1707", "1708", "1709")

$row = 16
foreach ($periodo in $periodos) {
    # Worker 1 (JHONNY JAVIER CARRILLO GARCIA) row for this periodo
    $ws.Range("B$row").Value = $ccDoc
    $ws.Range("C$row").Value = $doc1
    $ws.Range("D$row").Value = $name1
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $salarioBasico
    $row = $row + 1

    # Worker 2 (MAURICIO JAVIER TORRES ANGULO) row for the same periodo
    $ws.Range("B$row").Value = $ccDoc
    $ws.Range("C$row").Value = $doc2
    $ws.Range("D$row").Value = $name2
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $salarioBasico
    $row = $row + 1
}
